# Regenerate orders with updated distance/size codes.
# Distance codes: D51 -> D55, D64 -> D69, D80 -> D86
# Size code:      S30 -> S31  (S20 and S25 are unchanged)
#
# These substitutions apply everywhere the codes appear: the Condition
# names (e.g. Face17_D51_S30 -> Face17_D55_S31), the left/right filenames
# (e.g. Face17_D51_S30_l.png -> Face17_D55_S31_l.png,
# Fixation_D64_l.png -> Fixation_D69_l.png), and the standalone
# Distance/Size columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Use whole-cell-unaware substring replacement (xlPart) so codes embedded
# inside composite strings (Face17_D51_S30, Face17_D51_S30_l.png, ...) are
# updated along with the standalone Distance/Size values.
[void]$used.Replace("D51", "D55", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
[void]$used.Replace("D64", "D69", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
[void]$used.Replace("D80", "D86", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
[void]$used.Replace("S30", "S31", [Microsoft.Office.Interop.Excel.XlLookAt]::xlPart)
